$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Insert the new bold "Sujet 1 : ..." heading paragraph + blank paragraph
#    at the very start of the document.
# ---------------------------------------------------------------------------
$introXml = '<w:p ' + $wns + '>' +
    '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr>' +
    '<w:t>Sujet 1 : Les régimes totalitaires de l’URSS de Staline, de l’Allemagne nazie et de l’Italie fasciste ont-ils les mêmes caractéristiques ?</w:t>' +
    '</w:r>' +
    '</w:p>' +
    '<w:p ' + $wns + '/>'
$d.Range(0, 0).InsertXML($introXml)

# ---------------------------------------------------------------------------
# 2) Restructure the "Enfin, les régimes totalitaires ..." paragraph: move the
#    lastRenderedPageBreak to the start of the paragraph, and wrap "incité"
#    in a proofErr gramStart/gramEnd pair, splitting the trailing sentence
#    into its own runs.
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Enfin, les régimes totalitaires*") {
        $targetPara = $cand
        break
    }
}

$pStart = $targetPara.Range.Start
$pEndBefore = $targetPara.Range.End

$newBodyXml = '<w:p ' + $wns + '>' +
    '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Enfin, les régimes totalitaires ont des points communs et aussi des différences. La doctrine principale est la propagande, l’endoctrinement, a culte de la personnalité ainsi que le contrôle des masses sont des point communs assez fréquents. En revanche, l’idéologie, l’expansionnisme et le rôle de la femme vis-à-vis du régime est très différentes. Dans certains régimes, comme l’Allemagne, cette dernières est incitée à rester dans le foyer alors qu’en Italie et en URSS</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">, elle est </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>incité</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> à prêter main forte dans le monde du travail.</w:t></w:r>' +
    '</w:p>'
$d.Range($pStart, $pStart).InsertXML($newBodyXml)

$pEndAfter = $targetPara.Range.End
$insertedLen = $pEndAfter - $pEndBefore
$d.Range($pStart + $insertedLen, $pEndAfter).Delete()

# ---------------------------------------------------------------------------
# 3) Add header/footer parts (even/default/first) -- setting the default
#    header text creates header1/2/3.xml + footer1/2/3.xml + footnotes.xml +
#    endnotes.xml, matching what Word does the first time headers are turned
#    on for a document.
# ---------------------------------------------------------------------------
$section = $d.Sections.Item(1)
$section.PageSetup.DifferentFirstPageHeaderFooter = $true
$section.PageSetup.OddAndEvenPagesHeaderFooter = $true

# wdHeaderFooterPrimary = 1 -> header2.xml / footer2.xml ("default")
$primaryHeader = $section.Headers.Item(1)
$primaryHeader.Range.Text = "Henry Letellier T1`tHistoire`t27 10 2020"
